# Reorganize machinery template columns:
#  - Move "Operating Weight Range (kg)" from G to E (Canopy/Cab Version Weight shift right E->F, F->G)
#  - Move "Rated Power ISO14396 (kW)" from R to L (ISO9249/SAE J1349/EEC 80-1269/Number of
#    Cylinders/Bore x Stroke/Piston Displacement shift right L->M, M->N, N->O, O->P, P->Q, Q->R)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---

# E1:G1 block -> rotate right by one (Operating Weight Range moves to front)
$ws.Range("E1").Value = "Operating Weight Range (kg)"
$ws.Range("F1").Value = "Canopy Version Weight (kg)"
$ws.Range("G1").Value = "Cab Version Weight (kg)"

# L1:R1 block -> rotate right by one (Rated Power ISO14396 moves to front)
$ws.Range("L1").Value = "Rated Power ISO14396 (kW)"
$ws.Range("M1").Value = "Rated Power ISO9249 (kW)"
$ws.Range("N1").Value = "Rated Power SAE J1349 (kW)"
$ws.Range("O1").Value = "Rated Power EEC 80/1269 (kW)"
$ws.Range("P1").Value = "Number of Cylinders"
$ws.Range("Q1").Value = "Bore x Stroke (mm)"
$ws.Range("R1").Value = "Piston Displacement (L)"

# --- Data row (row 2) ---

# E2:G2 block -> rotate right by one
$ws.Range("E2").Value = 4000
$ws.Range("F2").Value = 3770
$ws.Range("G2").Value = 3940

# L2:R2 block -> rotate right by one
$ws.Range("L2").Value = 21.2
$ws.Range("M2").Value = 21.2
$ws.Range("N2").Value = 21.2
$ws.Range("O2").Value = 21.2
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = "88 x 90"
$ws.Range("R2").Value = 1.642
